$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update comments / precision-recall notes for the "adverb" related features ---
# Row 5: place adverbials
$ws.Range("D5").Value = "does not catch everything it should, but the reason seems to be the tagging, not the function"
$ws.Range("C5").Value = "there are some words that come to mind that could be added to Biber's list of place adverbials: apart, back, here, out, there (HM)"

# Row 6: time adverbials
$ws.Range("D6").Value = "does not catch everything it should, but the reason seems to be the tagging, not the function"
$ws.Range("C6").Value = "there are some words that come to mind that could be added to Biber's list of place adverbials: then, always (HM)"

# Row 36: adverbial subordinator of cause
$ws.Range("D36").Value = "does not catch everything it should, but the reason seems to be the tagging, not the function (""because"" tagged as IN)"
$ws.Range("C36").Value = "none"

# Row 37: adverbial subordinator of concession
$ws.Range("D37").Value = "does not catch everything it should, but the reason seems to be the tagging, not the function (""though"" tagged as IN)"
$ws.Range("C37").Value = "also included ""tho"" as common shortening of ""though"""

# --- Update sheet view: scroll position and selection ---
$ws.Activate()
$ws.Range("D37").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1

# --- Update workbook window placement (maximized-like window) ---
$wbWin = $wb.Windows.Item(1)
$wbWin.Left = -110
$wbWin.Top = -110
$wbWin.Width = 19420
$wbWin.Height = 10420
